# CDS Input file updates
# Replace the "ParticipantsTab" query (cell B2) with the new, updated Cypher
# query text (adds diagnosis/optional-match handling + sorted sample list),
# adjust row 2's height to fit the now-longer wrapped text, and move the
# active selection the way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newParticipantsQuery = @"
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE g.platform in ['Illumina NextSeq']
with p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN
coalesce(p.participant_id,'') as ``Participant ID``,
coalesce(s.study_name, '') as ``Study Name``,
coalesce(s.phs_accession,'') as ``Accession``,
coalesce(p.gender,'') as ``Gender``,
coalesce(apoc.text.join(samp, ','), '') as ``Samples``
ORDER BY p.participant_id LIMIT 100
"@

$ws.Range("B2").Value = $newParticipantsQuery

# The replacement query text wraps across more lines than the old one, so
# the row needs to grow to keep showing it in full.
$ws.Rows(2).RowHeight = 279

# Move the selection/scroll position the way it was left after the edit.
$ws.Range("B5").Select() | Out-Null
